$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values look numeric (e.g. "592.84"); a plain Value assignment
# that parses as a number gets auto-converted to a numeric cell by Excel.
# The source file stores these as plain text, so force each cell to text
# format before writing the new value.
$dCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D23", "D24", "D25", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D41", "D43", "D45", "D47", "D48", "D49", "D51")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.265.17"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "2.634.16"
$ws.Range("E3").Value = "  -2.79%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "592.84"
$ws.Range("E5").Value = "  -2.57%  "
$ws.Range("D6").Value = "166.76"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.536"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").Value = "2.633.61"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("D11").Value = "0.160"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").Value = "0.362"
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "5.24"
$ws.Range("E13").Value = "  -0.73%  "
$ws.Range("D14").Value = "27.73"
$ws.Range("E14").Value = "  -2.20%  "
$ws.Range("D15").Value = "3.113.42"
$ws.Range("E15").Value = "  -2.67%  "
$ws.Range("E16").Value = "  -2.71%  "
$ws.Range("D17").Value = "67.106.88"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "2.626.20"
$ws.Range("E18").Value = "  -2.70%  "
$ws.Range("D19").Value = "12.19"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("D20").Value = "8.12"
$ws.Range("E20").Value = "  +6.75%  "
$ws.Range("D21").Value = "360.68"
$ws.Range("E21").Value = "  -2.35%  "
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").Value = "4.69"
$ws.Range("E23").Value = "  -4.51%  "
$ws.Range("D24").Value = "10.89"
$ws.Range("E24").Value = "  +9.23%  "
$ws.Range("D25").Value = "1.96"
$ws.Range("E25").Value = "  -5.35%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "70.53"
$ws.Range("E27").Value = "  -3.32%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.0000102"
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "556.31"
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "7.96"
$ws.Range("E32").Value = "  -1.51%  "
$ws.Range("E33").Value = "  -2.74%  "
$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  -3.35%  "
$ws.Range("D35").Value = "0.135"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "1.52"
$ws.Range("E37").Value = "  -4.47%  "
$ws.Range("D38").Value = "157.65"
$ws.Range("E38").Value = "  -2.86%  "
$ws.Range("D39").Value = "19.21"
$ws.Range("E39").Value = "  -3.06%  "
$ws.Range("E40").Value = "  -2.62%  "
$ws.Range("D41").Value = "5.23"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("D43").Value = "17.92"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  -4.26%  "
$ws.Range("E46").Value = "  -1.36%  "
$ws.Range("D47").Value = "0.0₆0300"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "0.589"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("D49").Value = "152.49"
$ws.Range("E49").Value = "  -1.23%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "1.74"
$ws.Range("E51").Value = "  -1.09%  "

# Restore the default "Normal" style on the touched price cells so only
# the values changed (no stray number-format styling left behind).
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
